# Applies the "Added other modules in Ruby to selendroid" edit to the
# BatteryIndicator_JS.xlsx TestCases sheet:
#  - H2 (VT200-0251 validate text): validate4 gains two Iconposition checks
#  - G3/H3 (VT200-0254 script+validate): drop before/after screenshots in favour
#    of isIconDisplayed(true/false) checks
#  - G5/H5 (VT200-0256 script+validate): drop screenshot, add left Iconposition check
#  - G6/H6 (VT200-0257 script+validate): drop screenshot, add top Iconposition check
#  - Row 2 grows taller (180.75 -> 203.25) to fit the extra validate4 lines

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# H2
$ws.Cells.Item(2,8).Value = "validate1`n{`nvalidate_PageTitle=Compliance JS specs`n};`nvalidate2`n{`nvalidate_PageTitle=Battery JS Test`n};`nvalidate3`n{`nvalidate_Text_Exists=VT200-0251`n};`nvalidate4`n{`nvalidate_Screenshot=VT200_0251`nvalidate_Iconposition=batteryview_xpath,left,20`nvalidate_Iconposition=batteryview_xpath,top,40`n};"

# G3
$ws.Cells.Item(3,7).Value = "wait(3);`nvalidate1;`nlink_Click(battery_test_link);`nvalidate2;`nSelectTestToRun(VT200_0254_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nClickRunTest(runtest_bottom_xpath);`nwait(3);`nvalidate4;`nwait(10);`nvalidate5;"

# H3
$ws.Cells.Item(3,8).Value = "validate1`n{`nvalidate_PageTitle=Compliance JS specs`n};`nvalidate2`n{`nvalidate_PageTitle=Battery JS Test`n};`nvalidate3`n{`nvalidate_Text_Exists=VT200-0254`n};`nvalidate4`n{`nvalidate_isIconDisplayed=batteryview_xpath,true`n};`nvalidate5`n{`nvalidate_isIconDisplayed=batteryview_xpath,false`n};"

# G5
$ws.Cells.Item(5,7).Value = "wait(3);`nvalidate1;`nlink_Click(battery_test_link);`nvalidate2;`nSelectTestToRun(VT200_0256_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nClickRunTest(runtest_bottom_xpath);`nwait(3);`nvalidate4;"

# H5
$ws.Cells.Item(5,8).Value = "validate1`n{`nvalidate_PageTitle=Compliance JS specs`n};`nvalidate2`n{`nvalidate_PageTitle=Battery JS Test`n};`nvalidate3`n{`nvalidate_Text_Exists=VT200-0256`n};`nvalidate4`n{`nvalidate_Iconposition=batteryview_xpath,left,40`n};"

# G6
$ws.Cells.Item(6,7).Value = "wait(3);`nvalidate1;`nlink_Click(battery_test_link);`nvalidate2;`nSelectTestToRun(VT200_0257_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nClickRunTest(runtest_bottom_xpath);`nwait(3);`nvalidate4;`n"

# H6
$ws.Cells.Item(6,8).Value = "validate1`n{`nvalidate_PageTitle=Compliance JS specs`n};`nvalidate2`n{`nvalidate_PageTitle=Battery JS Test`n};`nvalidate3`n{`nvalidate_Text_Exists=VT200-0257`n};`nvalidate4`n{`nvalidate_Iconposition=batteryview_xpath,top,40`n};"

# Row 2 grew taller to fit the two extra Iconposition validate lines
$ws.Rows.Item(2).RowHeight = 203.25

